$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1), matching the style of the existing
# header cells (e.g. G1: bold, bordered, centered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new column's data rows (H2:H10) with 0, as plain numbers
# (no special style, matching the other numeric data columns).
$ws.Range("H2:H10").Value = 0
